$wb = $excel.ActiveWorkbook

$sheetRows = @{
    "ja_jp" = @(40, 44)
    "zh_cn" = @(40, 44)
    "en_gb" = @(44, 48)
    "ru_ru" = @(44, 48)
    "es_mx" = @(41, 45)
    "pt_br" = @(44, 48)
    "ar_ae" = @(44, 48)
    "de_de" = @(44, 48)
    "fr_fr" = @(44, 48)
    "en_ca" = @(43, 47)
    "fr_ca" = @(44, 48)
    "en_au" = @(41, 45)
    "en_in" = @(44, 48)
    "ko_kr" = @(44, 48)
    "en_sg" = @(44, 48)
    "no_no" = @(42, 46)
    "pl_pl" = @(41, 45)
    "pt_pt" = @(42, 46)
    "es_es" = @(42, 46)
    "sv_se" = @(42, 46)
    "fr_ch" = @(41, 45)
    "de_ch" = @(41, 45)
    "tr_tr" = @(42, 46)
    "ru_ua" = @(38, 42)
    "zh_hk" = @(42, 46)
    "en_hk" = @(42, 46)
    "zh_tw" = @(42, 46)
    "en_dz" = @(41, 45)
    "en_il" = @(38, 42)
    "nl_be" = @(44, 48)
    "da_dk" = @(44, 48)
    "fr_dz" = @(44, 48)
    "it_it" = @(44, 48)
    "en_ae" = @(44, 48)
    "nl_nl" = @(44, 48)
    "th_th" = @(44, 48)
    "de_at" = @(41, 45)
    "en_be" = @(41, 45)
    "fr_be" = @(41, 45)
}

foreach ($sheetName in $sheetRows.Keys) {
    $range = $sheetRows[$sheetName]
    $ws = $wb.Worksheets.Item($sheetName)
    $startRow = $range[0]
    $endRow = $range[1]
    $rangeAddr = "E" + $startRow + ":E" + $endRow
    $ws.Range($rangeAddr).Value = "YES"
}

Write-Host "Done updating sheets."